$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing D column values
$ws.Range("D2").Value = 0.001082428676132903
$ws.Range("D3").Value = 0.006346374577701655
$ws.Range("D4").Value = 0.02162982972700522
$ws.Range("D5").Value = 0.01781899850550332

# Update E5 value (tiny rounding change)
$ws.Range("E5").Value = 0.8661309802792048

# Add new row 6 - DWA (copy style from an existing labeled row cell, then set value)
$ws.Range("A2").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("A6").Value = "DWA"

$ws.Range("B6").Value = 100
$ws.Range("C6").Value = 0.9418132611637343
$ws.Range("D6").Value = 0.03023739526088455
$ws.Range("E6").Value = 0.9438686045747178
